$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1037.4166
$ws.Range("I15").Value = 1037.4166
$ws.Range("K15").Value = 3112.2498
$ws.Range("M15").Value = -2943.2498

$ws.Range("H17").Value = 923.275
$ws.Range("J17").Value = 923.275
$ws.Range("L17").Value = 2769.825
$ws.Range("N17").Value = -3105.825

$ws.Range("H40").Value = 12714.143
$ws.Range("I40").Value = 7333.1665
$ws.Range("K40").Value = 7333.1665
$ws.Range("M40").Value = -7158.1665

$ws.Range("H132").Value = 2658.0186
$ws.Range("I132").Value = 2692.5098
$ws.Range("K132").Value = 8077.529399999999
$ws.Range("M132").Value = -5547.529399999999

$ws.Range("H138").Value = 2925.4827
$ws.Range("I138").Value = 1575.2858
$ws.Range("J138").Value = 3691.8108
$ws.Range("K138").Value = 4725.857400000001
$ws.Range("L138").Value = 11075.4324
$ws.Range("M138").Value = 414.1425999999992
$ws.Range("N138").Value = -21355.4324

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6038.5415
$ws.Range("I32").Value = 2282.238
$ws.Range("K32").Value = 2282.238
$ws.Range("M32").Value = -1995.238

$ws.Range("H45").Value = 9849.379000000001
$ws.Range("I45").Value = 11840.381
$ws.Range("J45").Value = 4623
$ws.Range("K45").Value = 11840.381
$ws.Range("L45").Value = 4623
$ws.Range("M45").Value = -11463.381
$ws.Range("N45").Value = -5377

$ws.Range("H76").Value = 52500
$ws.Range("J76").Value = 52500
$ws.Range("L76").Value = 52500
$ws.Range("N76").Value = -53176

$ws.Range("H79").Value = 52500
$ws.Range("J79").Value = 52500
$ws.Range("L79").Value = 52500
$ws.Range("N79").Value = -54840

$ws.Range("H80").Value = 250000
$ws.Range("J80").Value = 250000
$ws.Range("L80").Value = 250000
$ws.Range("N80").Value = -251996

$ws.Range("H83").Value = 250000
$ws.Range("J83").Value = 250000
$ws.Range("L83").Value = 750000
$ws.Range("N83").Value = -759984

$ws.Range("H122").Value = 4160.5557
$ws.Range("I122").Value = 3920.8572
$ws.Range("K122").Value = 11762.5716
$ws.Range("M122").Value = -9312.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 89091.625
$ws.Range("J138").Value = 89091.625
$ws.Range("L138").Value = 89091.625
$ws.Range("N138").Value = -99371.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6266.4375
$ws.Range("I31").Value = 8161
$ws.Range("K31").Value = 8161
$ws.Range("M31").Value = -7866

$ws.Range("H34").Value = 6266.4375
$ws.Range("I34").Value = 8161
$ws.Range("K34").Value = 8161
$ws.Range("M34").Value = -7959

$ws.Range("H58").Value = 5118.8696
$ws.Range("I58").Value = 5346.421
$ws.Range("K58").Value = 5346.421
$ws.Range("M58").Value = -5143.421

$ws.Range("H107").Value = 1125.8334
$ws.Range("I107").Value = 1002
$ws.Range("J107").Value = 1373.5
$ws.Range("K107").Value = 1002
$ws.Range("L107").Value = 1373.5
$ws.Range("M107").Value = 918
$ws.Range("N107").Value = -5213.5

$ws.Range("H131").Value = 59994.25
$ws.Range("I131").Value = 39993
$ws.Range("K131").Value = 39993
$ws.Range("M131").Value = -34953

$ws.Range("H132").Value = 6645.614
$ws.Range("I132").Value = 6012.88
$ws.Range("J132").Value = 7478.1577
$ws.Range("K132").Value = 18038.64
$ws.Range("L132").Value = 22434.4731
$ws.Range("M132").Value = -15508.64
$ws.Range("N132").Value = -27494.4731

$ws.Range("H134").Value = 7753.574
$ws.Range("I134").Value = 6848.93
$ws.Range("J134").Value = 11289.909
$ws.Range("K134").Value = 20546.79
$ws.Range("L134").Value = 33869.727
$ws.Range("M134").Value = -18011.79
$ws.Range("N134").Value = -38939.727

$ws.Range("H135").Value = 80354
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 5118.8696
$ws.Range("I136").Value = 5346.421
$ws.Range("K136").Value = 16039.263
$ws.Range("M136").Value = -13489.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7278.976
$ws.Range("I56").Value = 7278.976
$ws.Range("K56").Value = 7278.976
$ws.Range("M56").Value = -6748.976

$ws.Range("H68").Value = 464.92856
$ws.Range("I68").Value = 350.91666
$ws.Range("J68").Value = 1149
$ws.Range("K68").Value = 1052.74998
$ws.Range("L68").Value = 3447
$ws.Range("M68").Value = -241.7499800000001
$ws.Range("N68").Value = -5069

$ws.Range("H71").Value = 464.92856
$ws.Range("I71").Value = 350.91666
$ws.Range("J71").Value = 1149
$ws.Range("K71").Value = 3158.24994
$ws.Range("L71").Value = 10341
$ws.Range("M71").Value = 897.7500600000003
$ws.Range("N71").Value = -18453

$ws.Range("H75").Value = 3429.625
$ws.Range("I75").Value = 3282.6667
$ws.Range("J75").Value = 3517.8
$ws.Range("K75").Value = 9848.000100000001
$ws.Range("L75").Value = 10553.4
$ws.Range("M75").Value = -8850.000100000001
$ws.Range("N75").Value = -12549.4

$ws.Range("H78").Value = 3429.625
$ws.Range("I78").Value = 3282.6667
$ws.Range("J78").Value = 3517.8
$ws.Range("K78").Value = 29544.0003
$ws.Range("L78").Value = 31660.2
$ws.Range("M78").Value = -24552.0003
$ws.Range("N78").Value = -41644.2

$ws.Range("H88").Value = 19999.666
$ws.Range("J88").Value = 19999.666
$ws.Range("L88").Value = 59998.99800000001
$ws.Range("N88").Value = -60854.99800000001

$ws.Range("H91").Value = 19999.666
$ws.Range("J91").Value = 19999.666
$ws.Range("L91").Value = 59998.99800000001
$ws.Range("N91").Value = -62962.99800000001

$ws.Range("H99").Value = 1099
$ws.Range("I99").Value = 1099
$ws.Range("K99").Value = 3297
$ws.Range("M99").Value = -1051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 69999
$ws.Range("J87").Value = 69999
$ws.Range("L87").Value = 69999
$ws.Range("N87").Value = -72495

$ws.Range("H90").Value = 69999
$ws.Range("J90").Value = 69999
$ws.Range("L90").Value = 209997
$ws.Range("N90").Value = -222477

$ws.Range("H126").Value = 6029.4375
$ws.Range("I126").Value = 5072
$ws.Range("J126").Value = 7260.4287
$ws.Range("K126").Value = 15216
$ws.Range("L126").Value = 21781.2861
$ws.Range("M126").Value = -12746
$ws.Range("N126").Value = -26721.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 760.1111
$ws.Range("I16").Value = 803.4375
$ws.Range("K16").Value = 803.4375
$ws.Range("M16").Value = -633.4375

$ws.Range("H18").Value = 34999.332
$ws.Range("I18").Value = 34999.332
$ws.Range("K18").Value = 34999.332
$ws.Range("M18").Value = -34827.332

$ws.Range("H20").Value = 4176332.8
$ws.Range("I20").Value = 8338666.5
$ws.Range("J20").Value = 13999.333
$ws.Range("K20").Value = 8338666.5
$ws.Range("L20").Value = 13999.333
$ws.Range("M20").Value = -8338440.5
$ws.Range("N20").Value = -14451.333

$ws.Range("H46").Value = 4816.933
$ws.Range("J46").Value = 2685
$ws.Range("L46").Value = 2685
$ws.Range("N46").Value = -3061

$ws.Range("H112").Value = 88748
$ws.Range("J112").Value = 88748
$ws.Range("L112").Value = 88748
$ws.Range("N112").Value = -91702

$ws.Range("H132").Value = 4498.1
$ws.Range("I132").Value = 1796.4
$ws.Range("J132").Value = 7199.8
$ws.Range("K132").Value = 5389.200000000001
$ws.Range("L132").Value = 21599.4
$ws.Range("M132").Value = -2859.200000000001
$ws.Range("N132").Value = -26659.4

$ws.Range("H141").Value = 88280.28999999999
$ws.Range("J141").Value = 88280.28999999999
$ws.Range("L141").Value = 88280.28999999999
$ws.Range("N141").Value = -98640.28999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 71662.336
$ws.Range("J95").Value = 71662.336
$ws.Range("L95").Value = 71662.336
$ws.Range("N95").Value = -77154.336

$ws.Range("H132").Value = 2583.9143
$ws.Range("I132").Value = 2754.4194
$ws.Range("J132").Value = 1262.5
$ws.Range("K132").Value = 8263.2582
$ws.Range("L132").Value = 3787.5
$ws.Range("M132").Value = -5733.2582
$ws.Range("N132").Value = -8847.5

$ws.Range("H136").Value = 3826.1316
$ws.Range("I136").Value = 3813.2666
$ws.Range("K136").Value = 11439.7998
$ws.Range("M136").Value = -8889.799800000001

$ws.Range("H137").Value = 88571.336
$ws.Range("J137").Value = 88571.336
$ws.Range("L137").Value = 88571.336
$ws.Range("N137").Value = -98771.336
